$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'45.350.35"
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = "'2.377.90"
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'317.30"
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = "'108.70"
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('D7').Value = "'0.639"
$ws.Range('E7').Value = '  +0.91%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('D10').Value = "'40.90"
$ws.Range('E10').Value = '  -4.10%  '
$ws.Range('D11').Value = "'0.0918"
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = "'8.53"
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = "'0.986"
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'15.58"
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = "'2.740.20"
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = "'2.373.72"
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = "'45.287.84"
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = "'15.69"
$ws.Range('E19').Value = '  +17.81%  '
$ws.Range('D20').Value = "'7.31"
$ws.Range('E20').Value = '  -3.69%  '
$ws.Range('D21').Value = "'0.0000106"
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').Value = "'3.68"
$ws.Range('E22').Value = '  +3.76%  '
$ws.Range('D23').Value = "'73.33"
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').Value = "'261.25"
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').Value = "'2.35"
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = "'7.63"
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').Value = "'11.25"
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'22.38"
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0960"
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').Value = "'37.26"
$ws.Range('E32').Value = '  -4.80%  '
$ws.Range('D33').Value = "'167.50"
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('D34').Value = "'2.88"
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('D36').Value = "'0.117"
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').Value = "'4.72"
$ws.Range('E37').Value = '  -4.17%  '
$ws.Range('D38').Value = "'4.05"
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').Value = "'1.92"
$ws.Range('E39').Value = '  +11.62%  '
$ws.Range('D40').Value = "'2.97"
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('D41').Value = "'0.0355"
$ws.Range('E41').Value = '  -2.44%  '
$ws.Range('D42').Value = "'97.92"
$ws.Range('E42').Value = '  -7.19%  '
$ws.Range('D43').Value = "'70.47"
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('D44').Value = "'13.05"
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').Value = "'0.229"
$ws.Range('E45').Value = '  -4.35%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = "'1.849.68"
$ws.Range('E46').Value = '  +12.56%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = "'6.00"
$ws.Range('E47').Value = '  +4.19%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').Value = "'1.00"
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = "'83.58"
$ws.Range('E49').Value = '  +6.35%  '
$ws.Range('D50').Value = "'112.23"
$ws.Range('E50').Value = '  -3.42%  '
$ws.Range('D51').Value = "'9.28"
$ws.Range('E51').Value = '  -0.03%  '
